$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: column letter, row number, new value
$updates = @(
    @('D', 2, '59.368.62'),
    @('E', 2, '  -5.65%  '),
    @('D', 3, '2.475.03'),
    @('E', 3, '  -7.75%  '),
    @('D', 4, '0.999'),
    @('E', 4, '  -0.11%  '),
    @('D', 5, '537.63'),
    @('E', 5, '  -2.79%  '),
    @('D', 6, '147.89'),
    @('E', 6, '  -6.70%  '),
    @('D', 7, '0.997'),
    @('E', 7, '  -0.18%  '),
    @('E', 8, '  -4.20%  '),
    @('D', 9, '2.474.88'),
    @('E', 9, '  -7.89%  '),
    @('D', 10, '0.0994'),
    @('E', 10, '  -5.96%  '),
    @('E', 11, '  -2.81%  '),
    @('D', 12, '5.33'),
    @('E', 12, '  -0.44%  '),
    @('E', 13, '  -4.35%  '),
    @('D', 14, '2.897.93'),
    @('E', 14, '  -8.11%  '),
    @('D', 15, '24.13'),
    @('E', 15, '  -8.13%  '),
    @('D', 16, '59.338.81'),
    @('E', 16, '  -5.53%  '),
    @('E', 17, '  -5.99%  '),
    @('D', 18, '2.518.47'),
    @('E', 18, '  -6.09%  '),
    @('D', 19, '11.17'),
    @('E', 19, '  -5.89%  '),
    @('D', 20, '4.34'),
    @('E', 20, '  -5.61%  '),
    @('D', 21, '323.75'),
    @('E', 21, '  -6.19%  '),
    @('D', 22, '0.968'),
    @('E', 22, '  -3.25%  '),
    @('D', 23, '5.73'),
    @('E', 23, '  -9.02%  '),
    @('D', 25, '60.62'),
    @('E', 25, '  -4.28%  '),
    @('D', 26, '0.162'),
    @('E', 26, '  -3.92%  '),
    @('D', 27, '0.979'),
    @('E', 27, '  -2.10%  '),
    @('D', 28, '7.73'),
    @('E', 28, '  -5.69%  '),
    @('E', 29, '  -5.78%  '),
    @('B', 30, 'PancakeSwap'),
    @('C', 30, 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'),
    @('D', 30, '1.82'),
    @('E', 30, '  -6.25%  '),
    @('B', 31, 'Fetch.AI'),
    @('C', 31, 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'),
    @('D', 31, '1.26'),
    @('E', 31, '  -9.16%  '),
    @('D', 32, '0.0₃0772'),
    @('E', 32, '  -10.20%  '),
    @('D', 33, '0.997'),
    @('E', 33, '  -0.16%  '),
    @('D', 34, '157.27'),
    @('E', 34, '  -5.22%  '),
    @('B', 35, 'NEARProtocol'),
    @('C', 35, 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'),
    @('D', 35, '4.59'),
    @('E', 35, '  -5.74%  '),
    @('B', 36, 'ImmutableX'),
    @('C', 36, 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @('D', 36, '1.38'),
    @('E', 36, '  -6.06%  '),
    @('D', 37, '18.38'),
    @('E', 37, '  -5.83%  '),
    @('D', 38, '1.75'),
    @('E', 38, '  -2.19%  '),
    @('D', 39, '5.97'),
    @('E', 39, '  -5.78%  '),
    @('D', 40, '318.73'),
    @('E', 40, '  -8.97%  '),
    @('D', 41, '36.77'),
    @('E', 41, '  -4.05%  '),
    @('D', 42, '0.841'),
    @('E', 42, '  -12.71%  '),
    @('D', 43, '3.71'),
    @('E', 43, '  -7.28%  '),
    @('D', 44, '0.998'),
    @('E', 44, '  -0.03%  '),
    @('E', 45, '  -2.69%  '),
    @('D', 46, '0.586'),
    @('E', 46, '  -5.08%  '),
    @('E', 47, '  -3.32%  '),
    @('D', 48, '0.0525'),
    @('E', 48, '  -6.55%  '),
    @('D', 49, '19.02'),
    @('E', 49, '  -8.99%  '),
    @('D', 50, '18.56'),
    @('E', 50, '  -8.76%  '),
    @('B', 51, 'Aave'),
    @('C', 51, 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'),
    @('D', 51, '121.98'),
    @('E', 51, '  -5.19%  ')
)

foreach ($u in $updates) {
    $col = $u[0]
    $row = $u[1]
    $val = $u[2]
    $cell = $ws.Range("$col$row")
    if ($col -eq "D") {
        # Column D holds price text that can look numeric (e.g. "0.999"),
        # so force text formatting before/after the write to keep it a string
        # like the source data, matching the original inlineStr cells.
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    } else {
        $cell.Value = $val
    }
}

Write-Output "Applied cryptos update: $($updates.Count) cells changed"